$d = $word.ActiveDocument

function Clear-StubRun($para) {
    # A freshly InsertParagraphAfter()-created blank paragraph serializes
    # with a leftover empty <w:r/>. Writing real text and then clearing it
    # again flushes that stub so the paragraph ends up with no run at all,
    # matching a genuinely empty paragraph.
    $rng = $para.Range
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = "X"
    $rng2 = $para.Range
    $rng2.MoveEnd(1, -1) | Out-Null
    $rng2.Text = ""
}

# ---------------------------------------------------------------------------
# 1) Merge the split runs (with proofErr wrappers) in the "1.1 CURRENT STATE
#    OF TECHNOLOGY IN OSA" body paragraph into a single run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(11)
$rng = $p.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = ""
$rng.Text = "The rapid growth of technology has contributed a lot to the continued progress of all classification of industry. However, some sectors today seem to be left behind in adapting the use of different methods to have progress on their part. Computerization is a control system that manages processes in the industrial workplace. It reduce human errors and processing time, thus it can boost productivity and result into a high quality of product produce. This can result in a system well integrated process that can perform much faster and more accurate than the manual system."

# ---------------------------------------------------------------------------
# 2) Merge the split runs in the "Web Based Room Reservation System is
#    appropriate..." paragraph into a single run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(13)
$rng = $p.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = ""
$rng.Text = "The Web Based Room Reservation System is appropriate for the use of students for faster and convenient room reservation for their own purpose. This system provide the user to reserve room for easy access base on the given schedule of the system."

# ---------------------------------------------------------------------------
# 3) Merge the split runs in the "The long term goal..." paragraph into a
#    single run.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(16)
$rng = $p.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = ""
$rng.Text = "The long term goal of the research is to develop an automated room reservation system for the students and its members of the available clubs and organizations in the school."

# ---------------------------------------------------------------------------
# 4) After "1.2.1 GENERAL OBJECTIVES" insert three new paragraphs: one
#    containing a lone tab, then two blank ones (matching the ind
#    firstLine=420 formatting already used by the surrounding paragraphs).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(17)
$p.Range.InsertParagraphAfter()
$tabPara = $d.Paragraphs.Item(18).Range
$tabPara.Collapse(1) | Out-Null
$tabPara.InsertAfter([char]9)

$p = $d.Paragraphs.Item(18)
$p.Range.InsertParagraphAfter()
Clear-StubRun $d.Paragraphs.Item(19)

$p = $d.Paragraphs.Item(19)
$p.Range.InsertParagraphAfter()
Clear-StubRun $d.Paragraphs.Item(20)

# ---------------------------------------------------------------------------
# 5) Rework the "1.3 SCOPE AND LIMITATIONS OF THE RESEARCH" heading
#    paragraph: it used to directly hold the heading text ("1" + ".3 SCOPE
#    ..."). It is now preceded by three new blank paragraphs (a
#    ListParagraph and two plain ones), and the heading's own runs get
#    merged into one.
# ---------------------------------------------------------------------------
$headingIdx = 0
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    if ($para.Range.Text -like "*1*.3 SCOPE AND LIMITATIONS OF THE RESEARCH*") {
        $headingIdx = $i
        break
    }
}

# Insert three new blank paragraphs after the heading (they inherit its
# plain "ind firstLine=420" formatting), then move the heading text down
# into the last of the three so the heading ends up after the new blanks.
$headingPara = $d.Paragraphs.Item($headingIdx)
$headingPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Item($headingIdx)
$headingPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Item($headingIdx)
$headingPara.Range.InsertParagraphAfter()

# Move the (merged) heading text into the 4th slot.
$target = $d.Paragraphs.Item($headingIdx + 3).Range
$target.MoveEnd(1, -1) | Out-Null
$target.Text = "1.3 SCOPE AND LIMITATIONS OF THE RESEARCH"

# Clear the original heading paragraph's own text.
$src = $d.Paragraphs.Item($headingIdx).Range
$src.MoveEnd(1, -1) | Out-Null
$src.Text = ""

# Turn the first (now-blank) paragraph into the ListParagraph/left=1140 one.
$listPara = $d.Paragraphs.Item($headingIdx)
$listPara.Range.ParagraphFormat.Style = "ListParagraph"
$listPara.Range.ParagraphFormat.LeftIndent = 57

# Clean up the stub runs on the two plain blank paragraphs in between.
Clear-StubRun $d.Paragraphs.Item($headingIdx + 1)
Clear-StubRun $d.Paragraphs.Item($headingIdx + 2)

$hi = $headingIdx + 3

# ---------------------------------------------------------------------------
# 6) Insert the brand-new paragraph (tab + the new scope/limitations body
#    text) right after the heading paragraph.
# ---------------------------------------------------------------------------
$headingPara = $d.Paragraphs.Item($hi)
$headingPara.Range.InsertParagraphAfter()
$bodyPara = $d.Paragraphs.Item($hi + 1)
$bodyRng = $bodyPara.Range
$bodyRng.Collapse(1) | Out-Null
$bodyRng.InsertAfter([char]9)

$bodyPara = $d.Paragraphs.Item($hi + 1)
$bodyRng = $bodyPara.Range
$bodyRng.MoveEnd(1, -1) | Out-Null
$bodyRng.Collapse(0) | Out-Null
$bodyText = "This study aims and introduces the Web Based Room Reservation System for Notre Dame of Dadiangas University. The web based room reservation system is consisting of three modules, the information, reports, and utilities. The information module let the users input ID number, name of organization or club, room number, and reason for reservation. The report module consists of information that needs to pass on to the PPO and notify if the room is available or not. The utilities module will be receive by the student. The student will fill up the utility form that includes number of hours, number of students, number of tables and chairs. The limitations of the said system are limited to the students and teachers of the Notre Dame of Dadiangas University. The system can only be use in the Office of Student Affairs."
$bodyRng.InsertAfter($bodyText)

# ---------------------------------------------------------------------------
# 7) Merge the split runs in the "1.5 DEFINITION OF TERMS" heading into a
#    single run.
# ---------------------------------------------------------------------------
$defPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*DEFINITION OF TERMS*") {
        $defPara = $para
        break
    }
}
$rng = $defPara.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = ""
$rng.Text = "1.5 DEFINITION OF TERMS   "
